$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.725.25"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "3.785.72"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'597.10"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").Value = "'169.33"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "3.785.58"
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.524"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("D11").Value = "'6.50"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "'0.453"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("E13").Value = "  +3.96%  "
$ws.Range("D14").Value = "'36.61"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "4.420.29"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("D16").Value = "3.782.42"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").Value = "'18.60"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "67.744.53"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "'7.19"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").Value = "'10.55"
$ws.Range("E21").Value = "  -6.59%  "
$ws.Range("D22").Value = "'469.31"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("E24").Value = "  -8.13%  "
$ws.Range("D25").Value = "'83.81"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").Value = "'2.20"
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").Value = "'2.91"
$ws.Range("D31").Value = "3.933.69"
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("D32").Value = "'7.63"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("D35").Value = "'9.14"
$ws.Range("E35").Value = "  -2.41%  "
$ws.Range("D36").Value = "3.746.37"
$ws.Range("E36").Value = "  -2.02%  "
$ws.Range("D37").Value = "'3.79"
$ws.Range("E37").Value = "  +1.92%  "
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "'0.311"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("E46").Value = "  -2.06%  "
$ws.Range("D47").Value = "'45.82"
$ws.Range("E47").Value = "  -2.77%  "
$ws.Range("D48").Value = "'396.16"
$ws.Range("E48").Value = "  -5.12%  "
$ws.Range("E49").Value = "  -8.15%  "
$ws.Range("D50").Value = "'140.40"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").Value = "'39.17"
$ws.Range("E51").Value = "  +3.11%  "
